$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Russia" column between Japan (K) and Saudi Arabia (was L, now M).
# This shifts the former L (Saudi Arabia) and M (USA) columns one to the right,
# carrying their header text and data along automatically.
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L1").Value = "Russia"

# Updated data (final render values for every row/column)
$ws.Range("B2").Value = 0.633130423220635
$ws.Range("C2").Value = 0.704067644137971
$ws.Range("D2").Value = 0.709905396398318
$ws.Range("E2").Value = 0.701126059691001
$ws.Range("F2").Value = 0.717787393269892
$ws.Range("G2").Value = 0.647161282776943
$ws.Range("H2").Value = 0.638620304460263
$ws.Range("I2").Value = 0.698210762411814
$ws.Range("J2").Value = 0.650902799584227
$ws.Range("K2").Value = 0.524491585481855
$ws.Range("L2").Value = 0.568973456184735
$ws.Range("M2").Value = 0.670833846908728
$ws.Range("N2").Value = 0.58055052974288

$ws.Range("B3").Value = 0.563498853260559
$ws.Range("C3").Value = 0.60475424371171
$ws.Range("D3").Value = 0.539474226981016
$ws.Range("E3").Value = 0.602478694572106
$ws.Range("F3").Value = 0.707890145357511
$ws.Range("G3").Value = 0.43983052336481
$ws.Range("H3").Value = 0.527172107887217
$ws.Range("I3").Value = 0.630421738216136
$ws.Range("J3").Value = 0.605304323854534
$ws.Range("K3").Value = 0.445837116675855
$ws.Range("L3").Value = 0.583955154990387
$ws.Range("M3").Value = 0.699951386105416
$ws.Range("N3").Value = 0.505635009871561

$ws.Range("B4").Value = 0.548888702456242
$ws.Range("C4").Value = 0.580080585349521
$ws.Range("D4").Value = 0.531485285723793
$ws.Range("E4").Value = 0.535520043103524
$ws.Range("F4").Value = 0.644766872360098
$ws.Range("G4").Value = 0.53660763349788
$ws.Range("H4").Value = 0.594661374773117
$ws.Range("I4").Value = 0.554344822951132
$ws.Range("J4").Value = 0.526043203777943
$ws.Range("K4").Value = 0.418778652482084
$ws.Range("L4").Value = 0.605453531287603
$ws.Range("M4").Value = 0.754595050954226
$ws.Range("N4").Value = 0.486093474963069

$ws.Range("A5").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("B5").Value = 0.492515087699993
$ws.Range("C5").Value = 0.523046603226999
$ws.Range("D5").Value = 0.472592378276595
$ws.Range("E5").Value = 0.427653406382078
$ws.Range("F5").Value = 0.605555864512875
$ws.Range("G5").Value = 0.521202526992176
$ws.Range("H5").Value = 0.525980483709787
$ws.Range("I5").Value = 0.553417384986919
$ws.Range("J5").Value = 0.512001347135006
$ws.Range("K5").Value = 0.358715278828469
$ws.Range("L5").Value = 0.5178833513058
$ws.Range("M5").Value = 0.703571225719179
$ws.Range("N5").Value = 0.457870290523818

$ws.Range("A6").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("B6").Value = 0.487666417243067
$ws.Range("C6").Value = 0.511419657456943
$ws.Range("D6").Value = 0.469893732599595
$ws.Range("E6").Value = 0.474588352501177
$ws.Range("F6").Value = 0.567811567214036
$ws.Range("G6").Value = 0.402783255465927
$ws.Range("H6").Value = 0.542201778861752
$ws.Range("I6").Value = 0.506093212662672
$ws.Range("J6").Value = 0.527599708110766
$ws.Range("K6").Value = 0.30423038698242
$ws.Range("L6").Value = 0.589896134442377
$ws.Range("M6").Value = 0.688471530144484
$ws.Range("N6").Value = 0.423889871127233

$ws.Range("A7").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("B7").Value = 0.486367361421124
$ws.Range("C7").Value = 0.576307051782642
$ws.Range("D7").Value = 0.542972705389941
$ws.Range("E7").Value = 0.562204961462218
$ws.Range("F7").Value = 0.688787338057347
$ws.Range("G7").Value = 0.475016652934936
$ws.Range("H7").Value = 0.466000641161711
$ws.Range("I7").Value = 0.570221750594586
$ws.Range("J7").Value = 0.507039845048516
$ws.Range("K7").Value = 0.408669702435669
$ws.Range("L7").Value = 0.34778801311712
$ws.Range("M7").Value = 0.530644213623694
$ws.Range("N7").Value = 0.425904969118064

$ws.Range("A8").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("B8").Value = 0.484425269846573
$ws.Range("C8").Value = 0.532496566334214
$ws.Range("D8").Value = 0.495088028905704
$ws.Range("E8").Value = 0.524921374189507
$ws.Range("F8").Value = 0.591941037014067
$ws.Range("G8").Value = 0.456819466801943
$ws.Range("H8").Value = 0.508200690995298
$ws.Range("I8").Value = 0.52300575337374
$ws.Range("J8").Value = 0.542266733123058
$ws.Range("K8").Value = 0.296876438769083
$ws.Range("L8").Value = 0.594355780533345
$ws.Range("M8").Value = 0.672464791241274
$ws.Range("N8").Value = 0.40374576344328

$ws.Range("A9").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("B9").Value = 0.472274429131213
$ws.Range("C9").Value = 0.540327123031372
$ws.Range("D9").Value = 0.561605121197964
$ws.Range("E9").Value = 0.484276172595847
$ws.Range("F9").Value = 0.58866177439978
$ws.Range("G9").Value = 0.422400589296013
$ws.Range("H9").Value = 0.524665035453288
$ws.Range("I9").Value = 0.51495362718769
$ws.Range("J9").Value = 0.54304619857183
$ws.Range("K9").Value = 0.280940372958279
$ws.Range("L9").Value = 0.456838794373566
$ws.Range("M9").Value = 0.603725863836496
$ws.Range("N9").Value = 0.429601311528009

$ws.Range("A10").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("B10").Value = 0.463687229299451
$ws.Range("C10").Value = 0.556554654633912
$ws.Range("D10").Value = 0.507372337420729
$ws.Range("E10").Value = 0.542203928429308
$ws.Range("F10").Value = 0.594223913846543
$ws.Range("G10").Value = 0.441717783939359
$ws.Range("H10").Value = 0.530543360898448
$ws.Range("I10").Value = 0.531455138277002
$ws.Range("J10").Value = 0.520220321402756
$ws.Range("K10").Value = 0.325524678081646
$ws.Range("L10").Value = 0.349016296615525
$ws.Range("M10").Value = 0.629350439518224
$ws.Range("N10").Value = 0.404921841368987

$ws.Range("B11").Value = 0.373752935747861
$ws.Range("C11").Value = 0.428037080634314
$ws.Range("D11").Value = 0.461656997855799
$ws.Range("E11").Value = 0.411908509173002
$ws.Range("F11").Value = 0.413269067027582
$ws.Range("G11").Value = 0.353766662905809
$ws.Range("H11").Value = 0.381200279280272
$ws.Range("I11").Value = 0.395162144580406
$ws.Range("J11").Value = 0.418574312646126
$ws.Range("K11").Value = 0.251461087147566
$ws.Range("L11").Value = 0.34371419848382
$ws.Range("M11").Value = 0.533333498726061
$ws.Range("N11").Value = 0.324748537110434

